$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion message (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.59 = 18080.14 pesos`n✅ 18080.14 pesos = 4.56 = 969.47 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate table values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 218
$wsTasas.Range("O10").Value = 3941.47
$wsTasas.Range("N12").Value = 3963
$wsTasas.Range("O12").Value = 212.5
